# ---------------------------------------------------------------------------
# CSAA Exam Docs Updated
#
# Adds a "DynamoDB best practices include:" section (one heading paragraph plus
# five bullet paragraphs) right after the existing
#   "- Only users and services can assume a role to take on permissions (not
#    groups)."
# paragraph, and gives that anchor paragraph a thin bottom rule (paragraph
# border) to visually separate it from the new section, matching the source
# diff exactly (including the proofErr grammar-check markers Word leaves around
# "time" / "possible").
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

# --- Locate the anchor paragraph -------------------------------------------
$anchorText = "Only users and services can assume a role"
$anchorIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "*$anchorText*") {
        $anchorIndex = $i
        break
    }
}
if ($anchorIndex -eq 0) {
    throw "Could not find the anchor paragraph containing: $anchorText"
}

# --- Give the anchor paragraph a bottom border ------------------------------
# <w:pBdr><w:bottom w:val="single" w:sz="6" w:space="1" w:color="auto"/></w:pBdr>
$anchorBorders = $d.Paragraphs($anchorIndex).Borders
$bottomBorder = $anchorBorders.Item(-3)   # wdBorderBottom
$bottomBorder.LineStyle = 1                # wdLineStyleSingle
$bottomBorder.LineWidth = 3                # -> w:sz="6" (eighths of a point)
$bottomBorder.ColorIndex = 0                # wdAuto -> w:color="auto"
$anchorBorders.DistanceFromBottom = 1       # -> w:space="1"

# --- Append the six new paragraphs, each as an exact OOXML fragment ---------
# Each fragment is inserted via Range.InsertXML into a freshly created empty
# paragraph, which keeps the paragraph/run formatting (and the <w:proofErr>
# grammar markers) identical to the authored source, rather than relying on
# inherited pPr/rPr formatting from typed text.
$newParagraphsXml = @(
    '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="F7F9FA"/><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Roboto" w:eastAsia="Times New Roman" w:hAnsi="Roboto" w:cs="Times New Roman"/><w:noProof w:val="0"/><w:color w:val="1C1D1F"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-GB" w:eastAsia="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Roboto" w:eastAsia="Times New Roman" w:hAnsi="Roboto" w:cs="Times New Roman"/><w:noProof w:val="0"/><w:color w:val="FF0000"/><w:sz w:val="32"/><w:szCs w:val="32"/><w:lang w:val="en-GB" w:eastAsia="en-GB"/></w:rPr><w:t xml:space="preserve">DynamoDB </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Roboto" w:eastAsia="Times New Roman" w:hAnsi="Roboto" w:cs="Times New Roman"/><w:noProof w:val="0"/><w:color w:val="1C1D1F"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-GB" w:eastAsia="en-GB"/></w:rPr><w:t>best practices include:</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>',
    '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="F7F9FA"/><w:spacing w:before="100" w:beforeAutospacing="1" w:after="100" w:afterAutospacing="1" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Roboto" w:eastAsia="Times New Roman" w:hAnsi="Roboto" w:cs="Times New Roman"/><w:noProof w:val="0"/><w:color w:val="1C1D1F"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-GB" w:eastAsia="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Roboto" w:eastAsia="Times New Roman" w:hAnsi="Roboto" w:cs="Times New Roman"/><w:noProof w:val="0"/><w:color w:val="1C1D1F"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-GB" w:eastAsia="en-GB"/></w:rPr><w:t>- Keep item sizes small.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>',
    '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="F7F9FA"/><w:spacing w:before="100" w:beforeAutospacing="1" w:after="100" w:afterAutospacing="1" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Roboto" w:eastAsia="Times New Roman" w:hAnsi="Roboto" w:cs="Times New Roman"/><w:noProof w:val="0"/><w:color w:val="1C1D1F"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-GB" w:eastAsia="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Roboto" w:eastAsia="Times New Roman" w:hAnsi="Roboto" w:cs="Times New Roman"/><w:noProof w:val="0"/><w:color w:val="1C1D1F"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-GB" w:eastAsia="en-GB"/></w:rPr><w:t>- If you are storing serial data in DynamoDB that will require actions based on data/</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Roboto" w:eastAsia="Times New Roman" w:hAnsi="Roboto" w:cs="Times New Roman"/><w:noProof w:val="0"/><w:color w:val="1C1D1F"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-GB" w:eastAsia="en-GB"/></w:rPr><w:t>time</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="Roboto" w:eastAsia="Times New Roman" w:hAnsi="Roboto" w:cs="Times New Roman"/><w:noProof w:val="0"/><w:color w:val="1C1D1F"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-GB" w:eastAsia="en-GB"/></w:rPr><w:t xml:space="preserve"> use separate tables for days, weeks, months.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>',
    '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="F7F9FA"/><w:spacing w:before="100" w:beforeAutospacing="1" w:after="100" w:afterAutospacing="1" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Roboto" w:eastAsia="Times New Roman" w:hAnsi="Roboto" w:cs="Times New Roman"/><w:noProof w:val="0"/><w:color w:val="1C1D1F"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-GB" w:eastAsia="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Roboto" w:eastAsia="Times New Roman" w:hAnsi="Roboto" w:cs="Times New Roman"/><w:noProof w:val="0"/><w:color w:val="1C1D1F"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-GB" w:eastAsia="en-GB"/></w:rPr><w:t>- Store more frequently and less frequently accessed data in separate tables.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>',
    '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="F7F9FA"/><w:spacing w:before="100" w:beforeAutospacing="1" w:after="100" w:afterAutospacing="1" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Roboto" w:eastAsia="Times New Roman" w:hAnsi="Roboto" w:cs="Times New Roman"/><w:noProof w:val="0"/><w:color w:val="1C1D1F"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-GB" w:eastAsia="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Roboto" w:eastAsia="Times New Roman" w:hAnsi="Roboto" w:cs="Times New Roman"/><w:noProof w:val="0"/><w:color w:val="1C1D1F"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-GB" w:eastAsia="en-GB"/></w:rPr><w:t xml:space="preserve">- If </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Roboto" w:eastAsia="Times New Roman" w:hAnsi="Roboto" w:cs="Times New Roman"/><w:noProof w:val="0"/><w:color w:val="1C1D1F"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-GB" w:eastAsia="en-GB"/></w:rPr><w:t>possible</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="Roboto" w:eastAsia="Times New Roman" w:hAnsi="Roboto" w:cs="Times New Roman"/><w:noProof w:val="0"/><w:color w:val="1C1D1F"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-GB" w:eastAsia="en-GB"/></w:rPr><w:t xml:space="preserve"> compress larger attribute values.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>',
    '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="F7F9FA"/><w:spacing w:before="100" w:beforeAutospacing="1" w:after="100" w:afterAutospacing="1" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Roboto" w:eastAsia="Times New Roman" w:hAnsi="Roboto" w:cs="Times New Roman"/><w:noProof w:val="0"/><w:color w:val="1C1D1F"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-GB" w:eastAsia="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Roboto" w:eastAsia="Times New Roman" w:hAnsi="Roboto" w:cs="Times New Roman"/><w:noProof w:val="0"/><w:color w:val="1C1D1F"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-GB" w:eastAsia="en-GB"/></w:rPr><w:t>- Store objects larger than 400KB in S3 and use pointers (S3 Object ID) in DynamoDB.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
)

$insertAfterIndex = $anchorIndex
foreach ($xmlFragment in $newParagraphsXml) {
    $d.Paragraphs($insertAfterIndex).Range.InsertParagraphAfter() | Out-Null
    $insertAfterIndex = $insertAfterIndex + 1
    $d.Paragraphs($insertAfterIndex).Range.InsertXML($xmlFragment) | Out-Null
}

Write-Output "Inserted $($newParagraphsXml.Count) DynamoDB best-practices paragraphs after paragraph $anchorIndex"
